$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (rows 6 and 7) - data set shrank from 6 to 4 rows
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Row 2: ECs / Efna4 / Epha5 / MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna4"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.180761
$ws.Range("H2").Value = 3.542283
$ws.Range("I2").Value = 0.6103536098015441
$ws.Range("J2").Value = 0.6103536098015441
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01136166666666667
$ws.Range("N2").Value = 0.034085
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.013415412895
$ws.Range("R2").Value = 0.120738716055
$ws.Range("S2").Value = 0.6103536098015441
$ws.Range("T2").Value = 0.6103536098015441

# Row 3: FAPs / Efna4 / Epha5 / MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna4"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4799286666666667
$ws.Range("H3").Value = 1.439786
$ws.Range("I3").Value = 0.2480825451952105
$ws.Range("J3").Value = 0.2480825451952105
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01136166666666667
$ws.Range("N3").Value = 0.034085
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.005452789534444444
$ws.Range("R3").Value = 0.04907510581
$ws.Range("S3").Value = 0.2480825451952105
$ws.Range("T3").Value = 0.2480825451952105

# Row 4: MuSCs / Efna4 / Epha5 / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Efna4"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2047916666666667
$ws.Range("H4").Value = 0.614375
$ws.Range("I4").Value = 0.1058599775968842
$ws.Range("J4").Value = 0.1058599775968842
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01136166666666667
$ws.Range("N4").Value = 0.034085
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.002326774652777778
$ws.Range("R4").Value = 0.020940971875
$ws.Range("S4").Value = 0.1058599775968842
$ws.Range("T4").Value = 0.1058599775968842

# Row 5: Resolving-Mac / Efna4 / Epha5 / MuSCs
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Efna4"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06907100000000001
$ws.Range("H5").Value = 0.207213
$ws.Range("I5").Value = 0.03570386740636119
$ws.Range("J5").Value = 0.03570386740636119
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01136166666666667
$ws.Range("N5").Value = 0.034085
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.0007847616783333333
$ws.Range("R5").Value = 0.007062855105
$ws.Range("S5").Value = 0.03570386740636119
$ws.Range("T5").Value = 0.03570386740636119
